$d = $word.ActiveDocument

# Replace author's surname, first name and patronymic in the Author-styled
# paragraph (each word is its own run, so MatchWholeWord keeps it precise).
$d.Content.Find.Execute("Беличева", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Апареев", 2)

$d.Content.Find.Execute("Дарья", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Дмитрий", 2)

$d.Content.Find.Execute("Михайловна", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Андреевич", 2)

# Fix verb gender agreement in the conclusion paragraph.
$d.Content.Find.Execute("я реализовала модель", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "я реализовал модель", 2)
